# Daily refresh of the Chicago carjacking-by-neighborhood-by-month workbook:
# advance the running "through <date>" month from Dec 07 to Dec 08, which
# updates the partial-December counts, nudges several other monthly totals
# that were revised, and re-sorts the neighborhood rows (column A) by rank,
# so rows 7-10 land on different neighborhoods than before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Sheet name and header label both carry the new "through" date.
$ws.Name = "Through 2022-12-08"
$ws.Range("B1").Value = "December 2022 (through December 08)"

# 2) Rows 7-10 are fully rewritten: the neighborhoods in those rows changed
#    (table re-sorted) and nearly every monthly count in them was revised,
#    so a single 4x97 array write is simplest and matches the final layout.
$rows7to10 = New-Object "object[,]" 4,97
# Row 7: Washington Heights
$rows7to10[0,0] = "Washington Heights"
$rows7to10[0,2] = 4
$rows7to10[0,3] = 3
$rows7to10[0,4] = 5
$rows7to10[0,6] = 3
$rows7to10[0,7] = 3
$rows7to10[0,9] = 3
$rows7to10[0,10] = 3
$rows7to10[0,11] = 3
$rows7to10[0,12] = 2
$rows7to10[0,13] = 1
$rows7to10[0,14] = 5
$rows7to10[0,15] = 4
$rows7to10[0,16] = 1
$rows7to10[0,17] = 1
$rows7to10[0,18] = 1
$rows7to10[0,19] = 2
$rows7to10[0,21] = 1
$rows7to10[0,23] = 1
$rows7to10[0,24] = 4
$rows7to10[0,25] = 1
$rows7to10[0,26] = 6
$rows7to10[0,27] = 6
$rows7to10[0,28] = 1
$rows7to10[0,29] = 1
$rows7to10[0,30] = 2
$rows7to10[0,31] = 2
$rows7to10[0,32] = 1
$rows7to10[0,33] = 2
$rows7to10[0,34] = 3
$rows7to10[0,36] = 1
$rows7to10[0,38] = 2
$rows7to10[0,39] = 3
$rows7to10[0,40] = 2
$rows7to10[0,41] = 1
$rows7to10[0,47] = 1
$rows7to10[0,48] = 2
$rows7to10[0,50] = 1
$rows7to10[0,53] = 1
$rows7to10[0,54] = 2
$rows7to10[0,56] = 2
$rows7to10[0,58] = 2
$rows7to10[0,59] = 1
$rows7to10[0,62] = 1
$rows7to10[0,63] = 3
$rows7to10[0,64] = 1
$rows7to10[0,65] = 1
$rows7to10[0,67] = 2
$rows7to10[0,68] = 3
$rows7to10[0,69] = 1
$rows7to10[0,70] = 1
$rows7to10[0,71] = 1
$rows7to10[0,72] = 2
$rows7to10[0,75] = 1
$rows7to10[0,76] = 2
$rows7to10[0,77] = 1
$rows7to10[0,78] = 1
$rows7to10[0,79] = 2
$rows7to10[0,80] = 2
$rows7to10[0,81] = 1
$rows7to10[0,83] = 2
$rows7to10[0,84] = 2
$rows7to10[0,85] = 1
$rows7to10[0,86] = 1
$rows7to10[0,90] = 1
$rows7to10[0,91] = 1
$rows7to10[0,92] = 1
# Row 8: Belmont Cragin
$rows7to10[1,0] = "Belmont Cragin"
$rows7to10[1,2] = 4
$rows7to10[1,3] = 1
$rows7to10[1,5] = 2
$rows7to10[1,6] = 1
$rows7to10[1,7] = 1
$rows7to10[1,9] = 3
$rows7to10[1,10] = 3
$rows7to10[1,11] = 2
$rows7to10[1,12] = 1
$rows7to10[1,16] = 1
$rows7to10[1,17] = 2
$rows7to10[1,18] = 1
$rows7to10[1,19] = 1
$rows7to10[1,20] = 3
$rows7to10[1,23] = 2
$rows7to10[1,26] = 2
$rows7to10[1,27] = 1
$rows7to10[1,29] = 2
$rows7to10[1,30] = 1
$rows7to10[1,31] = 2
$rows7to10[1,32] = 1
$rows7to10[1,35] = 1
$rows7to10[1,36] = 1
$rows7to10[1,37] = 1
$rows7to10[1,45] = 1
$rows7to10[1,46] = 2
$rows7to10[1,48] = 2
$rows7to10[1,52] = 1
$rows7to10[1,55] = 1
$rows7to10[1,56] = 1
$rows7to10[1,57] = 2
$rows7to10[1,59] = 1
$rows7to10[1,60] = 2
$rows7to10[1,63] = 1
$rows7to10[1,64] = 2
$rows7to10[1,65] = 1
$rows7to10[1,66] = 2
$rows7to10[1,67] = 2
$rows7to10[1,69] = 1
$rows7to10[1,70] = 1
$rows7to10[1,71] = 2
$rows7to10[1,72] = 2
$rows7to10[1,73] = 1
$rows7to10[1,74] = 1
$rows7to10[1,76] = 1
$rows7to10[1,80] = 1
$rows7to10[1,81] = 1
$rows7to10[1,83] = 1
$rows7to10[1,90] = 1
$rows7to10[1,92] = 1
$rows7to10[1,95] = 1
$rows7to10[1,96] = 1
# Row 9: South Shore
$rows7to10[2,0] = "South Shore"
$rows7to10[2,1] = 1
$rows7to10[2,2] = 4
$rows7to10[2,3] = 2
$rows7to10[2,4] = 6
$rows7to10[2,5] = 2
$rows7to10[2,6] = 7
$rows7to10[2,7] = 3
$rows7to10[2,8] = 5
$rows7to10[2,9] = 4
$rows7to10[2,10] = 2
$rows7to10[2,11] = 5
$rows7to10[2,12] = 6
$rows7to10[2,13] = 3
$rows7to10[2,14] = 4
$rows7to10[2,15] = 7
$rows7to10[2,16] = 3
$rows7to10[2,17] = 3
$rows7to10[2,18] = 6
$rows7to10[2,19] = 1
$rows7to10[2,20] = 5
$rows7to10[2,21] = 5
$rows7to10[2,22] = 4
$rows7to10[2,23] = 8
$rows7to10[2,24] = 9
$rows7to10[2,25] = 1
$rows7to10[2,26] = 9
$rows7to10[2,27] = 3
$rows7to10[2,28] = 4
$rows7to10[2,29] = 7
$rows7to10[2,30] = 5
$rows7to10[2,31] = 2
$rows7to10[2,33] = 3
$rows7to10[2,34] = 2
$rows7to10[2,35] = 1
$rows7to10[2,36] = 2
$rows7to10[2,40] = 4
$rows7to10[2,41] = 1
$rows7to10[2,42] = 5
$rows7to10[2,43] = 3
$rows7to10[2,44] = 4
$rows7to10[2,45] = 4
$rows7to10[2,47] = 1
$rows7to10[2,49] = 1
$rows7to10[2,51] = 3
$rows7to10[2,52] = 1
$rows7to10[2,54] = 3
$rows7to10[2,55] = 1
$rows7to10[2,57] = 1
$rows7to10[2,58] = 1
$rows7to10[2,61] = 1
$rows7to10[2,62] = 4
$rows7to10[2,63] = 6
$rows7to10[2,64] = 4
$rows7to10[2,65] = 5
$rows7to10[2,66] = 5
$rows7to10[2,67] = 9
$rows7to10[2,68] = 3
$rows7to10[2,69] = 2
$rows7to10[2,70] = 2
$rows7to10[2,71] = 3
$rows7to10[2,72] = 4
$rows7to10[2,73] = 2
$rows7to10[2,74] = 5
$rows7to10[2,75] = 4
$rows7to10[2,77] = 3
$rows7to10[2,79] = 1
$rows7to10[2,80] = 1
$rows7to10[2,83] = 1
$rows7to10[2,84] = 1
$rows7to10[2,87] = 3
$rows7to10[2,89] = 1
$rows7to10[2,90] = 1
$rows7to10[2,91] = 1
$rows7to10[2,92] = 2
$rows7to10[2,93] = 1
$rows7to10[2,94] = 1
$rows7to10[2,96] = 1
# Row 10: Grand Boulevard
$rows7to10[3,0] = "Grand Boulevard"
$rows7to10[3,2] = 3
$rows7to10[3,3] = 2
$rows7to10[3,4] = 3
$rows7to10[3,5] = 1
$rows7to10[3,6] = 2
$rows7to10[3,7] = 1
$rows7to10[3,8] = 4
$rows7to10[3,9] = 2
$rows7to10[3,10] = 1
$rows7to10[3,11] = 1
$rows7to10[3,12] = 2
$rows7to10[3,13] = 3
$rows7to10[3,14] = 5
$rows7to10[3,15] = 5
$rows7to10[3,16] = 3
$rows7to10[3,17] = 3
$rows7to10[3,18] = 2
$rows7to10[3,19] = 1
$rows7to10[3,20] = 1
$rows7to10[3,22] = 2
$rows7to10[3,23] = 2
$rows7to10[3,24] = 9
$rows7to10[3,26] = 8
$rows7to10[3,27] = 2
$rows7to10[3,28] = 3
$rows7to10[3,29] = 3
$rows7to10[3,30] = 1
$rows7to10[3,31] = 3
$rows7to10[3,32] = 2
$rows7to10[3,34] = 3
$rows7to10[3,36] = 1
$rows7to10[3,40] = 2
$rows7to10[3,42] = 1
$rows7to10[3,44] = 1
$rows7to10[3,45] = 2
$rows7to10[3,57] = 1
$rows7to10[3,60] = 2
$rows7to10[3,62] = 1
$rows7to10[3,63] = 1
$rows7to10[3,64] = 4
$rows7to10[3,65] = 1
$rows7to10[3,70] = 1
$rows7to10[3,74] = 4
$rows7to10[3,75] = 3
$rows7to10[3,76] = 1
$rows7to10[3,78] = 1
$rows7to10[3,79] = 1
$rows7to10[3,80] = 1
$rows7to10[3,89] = 2
$rows7to10[3,90] = 1
$rows7to10[3,91] = 1
$rows7to10[3,94] = 1
$rows7to10[3,95] = 1
$ws.Range("A7:CS10").Value = $rows7to10

# 3) Remaining scattered single-cell revisions elsewhere in the table.
$ws.Range("BJ3").Value = 3
$ws.Range("CH3").Value = 1
$ws.Range("B4").Value = 1
$ws.Range("Z4").Value = 1
$ws.Range("BV5").Value = 3
$ws.Range("N6").Value = 3
$ws.Range("B12").Value = 2
$ws.Range("B14").Value = 2
$ws.Range("Z14").Value = 4
$ws.Range("N15").Value = 4
$ws.Range("N16").Value = 2
$ws.Range("Z20").Value = 6
$ws.Range("AL20").Value = 1
$ws.Range("BJ21").Value = 2
$ws.Range("AX28").Value = 1
$ws.Range("BJ28").Value = 4
$ws.Range("B30").Value = 2
$ws.Range("B35").Value = 1
$ws.Range("N35").Value = 1
$ws.Range("BJ37").Value = 2
$ws.Range("AX40").Value = 1
$ws.Range("AX43").Value = 1
$ws.Range("B47").Value = 1
$ws.Range("B49").Value = 1
$ws.Range("N61").Value = 2
$ws.Range("N64").Value = 4
$ws.Range("CH64").Value = 1
$ws.Range("N83").Value = 2
$ws.Range("AL87").Value = 1
$ws.Range("B88").Value = 1
